$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5196
$ws.Range("J40").Value = 4134.3335
$ws.Range("L40").Value = 4134.3335
$ws.Range("N40").Value = -4484.3335

$ws.Range("H64").Value = 5312.2
$ws.Range("I64").Value = 2943.7778
$ws.Range("J64").Value = 7250
$ws.Range("K64").Value = 2943.7778
$ws.Range("L64").Value = 7250
$ws.Range("M64").Value = -2695.7778
$ws.Range("N64").Value = -7746

$ws.Range("H67").Value = 5312.2
$ws.Range("I67").Value = 2943.7778
$ws.Range("J67").Value = 7250
$ws.Range("K67").Value = 2943.7778
$ws.Range("L67").Value = 7250
$ws.Range("M67").Value = -2085.7778
$ws.Range("N67").Value = -8966

$ws.Range("H86").Value = 1717.8649
$ws.Range("I86").Value = 1933.1904
$ws.Range("J86").Value = 1435.25
$ws.Range("K86").Value = 1933.1904
$ws.Range("L86").Value = 1435.25
$ws.Range("M86").Value = -810.1904
$ws.Range("N86").Value = -3681.25

$ws.Range("H88").Value = 2647.8708
$ws.Range("J88").Value = 2775.1482
$ws.Range("L88").Value = 2775.1482
$ws.Range("N88").Value = -3587.1482

$ws.Range("H89").Value = 1717.8649
$ws.Range("I89").Value = 1933.1904
$ws.Range("J89").Value = 1435.25
$ws.Range("K89").Value = 9665.951999999999
$ws.Range("L89").Value = 7176.25
$ws.Range("M89").Value = -4049.951999999999
$ws.Range("N89").Value = -18408.25

$ws.Range("H91").Value = 2647.8708
$ws.Range("J91").Value = 2775.1482
$ws.Range("L91").Value = 2775.1482
$ws.Range("N91").Value = -5583.1482

$ws.Range("H96").Value = 742
$ws.Range("I96").Value = 645
$ws.Range("J96").Value = 887.5
$ws.Range("K96").Value = 1935
$ws.Range("L96").Value = 2662.5
$ws.Range("M96").Value = -562
$ws.Range("N96").Value = -5408.5

$ws.Range("H137").Value = 48824.42
$ws.Range("I137").Value = 105891.414
$ws.Range("J137").Value = 2627.3333
$ws.Range("K137").Value = 317674.242
$ws.Range("L137").Value = 7881.999899999999
$ws.Range("M137").Value = -315124.242
$ws.Range("N137").Value = -12981.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4822.5
$ws.Range("I32").Value = 3519.6047
$ws.Range("J32").Value = 7258.3477
$ws.Range("K32").Value = 3519.6047
$ws.Range("L32").Value = 7258.3477
$ws.Range("M32").Value = -3232.6047
$ws.Range("N32").Value = -7832.3477

$ws.Range("H45").Value = 7574419
$ws.Range("J45").Value = 5536.5713
$ws.Range("L45").Value = 5536.5713
$ws.Range("N45").Value = -6290.5713

$ws.Range("H55").Value = 47999
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 47999
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 47999
$ws.Range("N55").Value = -48629
$ws.Range("M55").ClearContents()

$ws.Range("H61").Value = 2368.125
$ws.Range("J61").Value = 3706.75
$ws.Range("L61").Value = 3706.75
$ws.Range("N61").Value = -4130.75

$ws.Range("H132").Value = 2447.5715
$ws.Range("I132").Value = 1468.4166
$ws.Range("J132").Value = 3753.111
$ws.Range("K132").Value = 4405.2498
$ws.Range("L132").Value = 11259.333
$ws.Range("M132").Value = -1875.2498
$ws.Range("N132").Value = -16319.333

$ws.Range("H136").Value = 2368.125
$ws.Range("J136").Value = 3706.75
$ws.Range("L136").Value = 11120.25
$ws.Range("N136").Value = -16220.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3130746.2
$ws.Range("I86").Value = 4006590.8
$ws.Range("J86").Value = 2729.7144
$ws.Range("K86").Value = 4006590.8
$ws.Range("L86").Value = 2729.7144
$ws.Range("M86").Value = -4005467.8
$ws.Range("N86").Value = -4975.7144

$ws.Range("H89").Value = 3130746.2
$ws.Range("I89").Value = 4006590.8
$ws.Range("J89").Value = 2729.7144
$ws.Range("K89").Value = 20032954
$ws.Range("L89").Value = 13648.572
$ws.Range("M89").Value = -20027338
$ws.Range("N89").Value = -24880.572

$ws.Range("H99").Value = 8993961
$ws.Range("I99").Value = 13080325
$ws.Range("J99").Value = 3959
$ws.Range("K99").Value = 13080325
$ws.Range("L99").Value = 3959
$ws.Range("M99").Value = -13078827
$ws.Range("N99").Value = -6955

$ws.Range("H134").Value = 3391.9211
$ws.Range("I134").Value = 1732.4193
$ws.Range("K134").Value = 5197.257900000001
$ws.Range("M134").Value = -2662.257900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4379.4
$ws.Range("I99").Value = 4333.3335
$ws.Range("J99").Value = 4399.143
$ws.Range("K99").Value = 4333.3335
$ws.Range("L99").Value = 4399.143
$ws.Range("M99").Value = -2835.3335
$ws.Range("N99").Value = -7395.143

$ws.Range("H126").Value = 4379.4
$ws.Range("I126").Value = 4333.3335
$ws.Range("J126").Value = 4399.143
$ws.Range("K126").Value = 13000.0005
$ws.Range("L126").Value = 13197.429
$ws.Range("M126").Value = -10530.0005
$ws.Range("N126").Value = -18137.429

$ws.Range("H132").Value = 88466.95
$ws.Range("I132").Value = 68609.336
$ws.Range("J132").Value = 131019
$ws.Range("K132").Value = 205828.008
$ws.Range("L132").Value = 393057
$ws.Range("M132").Value = -203298.008
$ws.Range("N132").Value = -398117

$ws.Range("H134").Value = 3957
$ws.Range("I134").Value = 3057.0667
$ws.Range("J134").Value = 5885.4287
$ws.Range("K134").Value = 9171.2001
$ws.Range("L134").Value = 17656.2861
$ws.Range("M134").Value = -6636.2001
$ws.Range("N134").Value = -22726.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 63777.75
$ws.Range("I5").Value = 1044.7
$ws.Range("J5").Value = 168332.83
$ws.Range("K5").Value = 3134.1
$ws.Range("L5").Value = 504998.49
$ws.Range("M5").Value = -3022.1
$ws.Range("N5").Value = -505222.49

$ws.Range("H7").Value = 1229.6923
$ws.Range("I7").Value = 1951.3334
$ws.Range("J7").Value = 611.1429000000001
$ws.Range("K7").Value = 5854.0002
$ws.Range("L7").Value = 1833.4287
$ws.Range("M7").Value = -5742.0002
$ws.Range("N7").Value = -2057.4287

$ws.Range("H92").Value = 653.55554
$ws.Range("I92").Value = 537.6667
$ws.Range("J92").Value = 711.5
$ws.Range("K92").Value = 1613.0001
$ws.Range("L92").Value = 2134.5
$ws.Range("M92").Value = -365.0001
$ws.Range("N92").Value = -4630.5

$ws.Range("H135").Value = 63777.75
$ws.Range("I135").Value = 1044.7
$ws.Range("J135").Value = 168332.83
$ws.Range("K135").Value = 9402.300000000001
$ws.Range("L135").Value = 1514995.47
$ws.Range("M135").Value = -6867.300000000001
$ws.Range("N135").Value = -1520065.47

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 24333.334
$ws.Range("I44").Value = 24333.334
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 24333.334
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -23737.334
$ws.Range("N44").ClearContents()

$ws.Range("H70").Value = 13340546
$ws.Range("I70").Value = 22229922
$ws.Range("J70").Value = 6482.3335
$ws.Range("K70").Value = 22229922
$ws.Range("L70").Value = 6482.3335
$ws.Range("M70").Value = -22229652
$ws.Range("N70").Value = -7022.3335

$ws.Range("H73").Value = 13340546
$ws.Range("I73").Value = 22229922
$ws.Range("J73").Value = 6482.3335
$ws.Range("K73").Value = 22229922
$ws.Range("L73").Value = 6482.3335
$ws.Range("M73").Value = -22228986
$ws.Range("N73").Value = -8354.333500000001

$ws.Range("H122").Value = 373019.8
$ws.Range("I122").Value = 567225.7
$ws.Range("J122").Value = 4028.7
$ws.Range("K122").Value = 1701677.1
$ws.Range("L122").Value = 12086.1
$ws.Range("M122").Value = -1699227.1
$ws.Range("N122").Value = -16986.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 880.71875
$ws.Range("I16").Value = 636.9286
$ws.Range("J16").Value = 2587.25
$ws.Range("K16").Value = 636.9286
$ws.Range("L16").Value = 2587.25
$ws.Range("M16").Value = -466.9286
$ws.Range("N16").Value = -2927.25

$ws.Range("H46").Value = 4444.278
$ws.Range("J46").Value = 5007.4614
$ws.Range("L46").Value = 5007.4614
$ws.Range("N46").Value = -5383.4614

$ws.Range("H61").Value = 6539183
$ws.Range("I61").Value = 9263009
$ws.Range("J61").Value = 1999.6
$ws.Range("K61").Value = 9263009
$ws.Range("L61").Value = 1999.6
$ws.Range("M61").Value = -9262807
$ws.Range("N61").Value = -2403.6

$ws.Range("H113").Value = 6539183
$ws.Range("I113").Value = 9263009
$ws.Range("J113").Value = 1999.6
$ws.Range("K113").Value = 9263009
$ws.Range("L113").Value = 1999.6
$ws.Range("M113").Value = -9260839
$ws.Range("N113").Value = -6339.6

$ws.Range("H132").Value = 3857.7358
$ws.Range("I132").Value = 3113.0889
$ws.Range("J132").Value = 8046.375
$ws.Range("K132").Value = 9339.2667
$ws.Range("L132").Value = 24139.125
$ws.Range("M132").Value = -6809.2667
$ws.Range("N132").Value = -29199.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 350
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H132").Value = 126365736
$ws.Range("I132").Value = 250003500
$ws.Range("K132").Value = 750010500
$ws.Range("M132").Value = -750007970
